# Populate the new "Category.field" value (column Y) on the Data sheet
# for every data row (2-11) with "object_annotation_category".
# This also implicitly registers a new shared string in sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 25).Value2 = "object_annotation_category"
}
